$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - use leading apostrophe so the numeric-looking
# text is stored as text (matching the source data's inline-string cells)
# instead of being auto-converted to a real number by Excel.
$ws.Range("D2").Value = "'235.58"
$ws.Range("D3").Value = "'22.25"
$ws.Range("D4").Value = "'5.424"
$ws.Range("D5").Value = "'0.05643"
$ws.Range("D6").Value = "'6.485"
$ws.Range("D8").Value = "'1.067"
$ws.Range("D9").Value = "'0.7870"
$ws.Range("D10").Value = "'0.1398"
$ws.Range("D11").Value = "'0.07334"
$ws.Range("D12").Value = "'0.03210"
$ws.Range("D13").Value = "'0.02975"
$ws.Range("D14").Value = "'0.09251"
$ws.Range("D15").Value = "'0.001660"
$ws.Range("D16").Value = "'3.263"
$ws.Range("D17").Value = "'0.04756"
$ws.Range("D18").Value = "'0.0005757"
$ws.Range("D19").Value = "'0.006213"
$ws.Range("D20").Value = "'0.005097"
$ws.Range("D23").Value = "'3.856"
$ws.Range("D24").Value = "'2.147"
$ws.Range("D40").Value = "'0.04099"
$ws.Range("D41").Value = "'0.006971"
$ws.Range("D44").Value = "'0.009922"
$ws.Range("D45").Value = "'0.00005429"
$ws.Range("D47").Value = "'0.6760"
$ws.Range("D48").Value = "'0.03871"

# Volume(1h) (column E) text-only updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

# Rows 42/43: BKEXToken and CEJI swapped places in the ranking
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003504"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1038"
$ws.Range("E43").Value = "42BKEXTokenBKK"
